$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$before = $ws.Cells.Item(7,4).Value()
Write-Host ("Before D7: {0}" -f $before)
$ws.Columns("D").Insert()
$afterD = $ws.Cells.Item(7,4).Value()
$afterE = $ws.Cells.Item(7,5).Value()
Write-Host ("After insert D7: {0}" -f $afterD)
Write-Host ("After insert E7: {0}" -f $afterE)
